$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Vtn", "Itgb3", "ECs", 2, 0.6666666666666666, 6.597131, 19.791393, 0.1209543635982448, 0.1209543635982448, 3, 1, 2.481489333333333, 7.444467999999999, 0.2345069082418988, 0.2345069082418987, 16.37071020710266, 147.336391863924, 0.02836463384579085, 0.02836463384579084),
    @("ECs", "Vtn", "Itgb3", "FAPs", 2, 0.6666666666666666, 6.597131, 19.791393, 0.1209543635982448, 0.1209543635982448, 3, 1, 7.245227, 21.735681, 0.6846919551326144, 0.6846919551326142, 47.797711643737, 430.179404793633, 0.08281647969390334, 0.08281647969390331),
    @("ECs", "Vtn", "Itgb3", "M2", 2, 0.6666666666666666, 6.597131, 19.791393, 0.1209543635982448, 0.1209543635982448, 2, 0.6666666666666666, 0.2001876666666667, 0.600563, 0.01891823194544989, 0.01891823194544989, 1.320664261584333, 11.885978354259, 0.002288242705365876, 0.002288242705365875),
    @("ECs", "Vtn", "Itgb3", "sCs", 2, 0.6666666666666666, 6.597131, 19.791393, 0.1209543635982448, 0.1209543635982448, 3, 1, 0.6548283333333333, 1.964485, 0.06188290468003712, 0.06188290468003711, 4.319988297511667, 38.879894677605, 0.007485007353184732, 0.007485007353184731),
    @("FAPs", "Vtn", "Itgb3", "ECs", 3, 1, 21.05317333333333, 63.15952, 0.3859970617919927, 0.3859970617919927, 3, 1, 2.481489333333333, 7.444467999999999, 0.2345069082418988, 0.2345069082418987, 52.24322505948443, 470.1890255353599, 0.09051897755129737, 0.09051897755129734),
    @("FAPs", "Vtn", "Itgb3", "FAPs", 3, 1, 21.05317333333333, 63.15952, 0.3859970617919927, 0.3859970617919927, 3, 1, 7.245227, 21.735681, 0.6846919551326144, 0.6846919551326142, 152.5350198703467, 1372.81517883312, 0.264289082913804, 0.264289082913804),
    @("FAPs", "Vtn", "Itgb3", "M2", 3, 1, 21.05317333333333, 63.15952, 0.3859970617919927, 0.3859970617919927, 2, 0.6666666666666666, 0.2001876666666667, 0.600563, 0.01891823194544989, 0.01891823194544989, 4.214585645528889, 37.93127080976, 0.007302381945243073, 0.007302381945243071),
    @("FAPs", "Vtn", "Itgb3", "sCs", 3, 1, 21.05317333333333, 63.15952, 0.3859970617919927, 0.3859970617919927, 3, 1, 0.6548283333333333, 1.964485, 0.06188290468003712, 0.06188290468003711, 13.78621440524445, 124.0759296472, 0.02388661938164828, 0.02388661938164828),
    @("sCs", "Vtn", "Itgb3", "ECs", 3, 1, 26.89201066666667, 80.676032, 0.4930485746097625, 0.4930485746097625, 3, 1, 2.481489333333333, 7.444467999999999, 0.2345069082418988, 0.2345069082418987, 66.73223762121954, 600.5901385909759, 0.1156232968448105, 0.1156232968448105),
    @("sCs", "Vtn", "Itgb3", "FAPs", 3, 1, 26.89201066666667, 80.676032, 0.4930485746097625, 0.4930485746097625, 3, 1, 7.245227, 21.735681, 0.6846919551326144, 0.6846919551326142, 194.8387217664213, 1753.548495897792, 0.337586392524907, 0.3375863925249069),
    @("sCs", "Vtn", "Itgb3", "M2", 3, 1, 26.89201066666667, 80.676032, 0.4930485746097625, 0.4930485746097625, 2, 0.6666666666666666, 0.2001876666666667, 0.600563, 0.01891823194544989, 0.01891823194544989, 5.383448867335111, 48.451039806016, 0.009327607294840943, 0.009327607294840942),
    @("sCs", "Vtn", "Itgb3", "sCs", 3, 1, 26.89201066666667, 80.676032, 0.4930485746097625, 0.4930485746097625, 3, 1, 0.6548283333333333, 1.964485, 0.06188290468003712, 0.06188290468003711, 17.60965052483555, 158.48685472352, 0.0305112779452041, 0.0305112779452041)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowData = $data[$r]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $rowData[$c]
    }
}
